$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheetPr/outlinePr defaults used on the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1

# Match the pageMargins used on the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in).
$psForecast = $wsForecast.PageSetup
$psForecast.LeftMargin = 54
$psForecast.RightMargin = 54
$psForecast.TopMargin = 72
$psForecast.BottomMargin = 72
$psForecast.HeaderMargin = 36
$psForecast.FooterMargin = 36

# Copy formatting (bold/border/alignment) from an existing header row onto the
# new header row so the same cell style is reused, then set the new header text.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:D1"))
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the date-number-format style from column A of the weekly sheet down the
# whole new date column so the cells reuse the same numFmt style.
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A11"))

# --- Populate data rows ---
$dates = @(44934.99999999999, 44983.99999999999, 44990.99999999999, 44997.99999999999, 45004.99999999999, 45011.99999999999, 45018.99999999999, 45025.99999999999, 45032.99999999999, 45039.99999999999)
$lowers = @(19.9999999766618, 19.9999999750047, 19.99999997142289, 19.99999996886211, 19.99999995918274, 19.99999994904554, 19.99999991945164, 19.99999989409661, 19.99999985702339, 19.99999982161052)
$uppers = @(20.0000000253745, 20.00000002526136, 20.00000002801901, 20.0000000293998, 20.00000003746174, 20.00000004671624, 20.00000006606727, 20.00000010191842, 20.00000014856125, 20.00000018752734)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = 20
    $wsForecast.Cells.Item($row, 3).Value = $lowers[$i]
    $wsForecast.Cells.Item($row, 4).Value = $uppers[$i]
}

Write-Host "done"
